# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the "Balmung_Profits" workbook per the commit diff.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 5561444.5
$ws.Range("I9").Value = 8333685
$ws.Range("J9").Value = 16963.6
$ws.Range("K9").Value = 8333685
$ws.Range("L9").Value = 16963.6
$ws.Range("M9").Value = -8333516
$ws.Range("N9").Value = -17301.6
$ws.Range("H42").Value = 801.3333
$ws.Range("I42").Value = 77
$ws.Range("K42").Value = 231
$ws.Range("M42").Value = -1
$ws.Range("H64").Value = 4179.3
$ws.Range("I64").Value = 4113.7144
$ws.Range("K64").Value = 4113.7144
$ws.Range("M64").Value = -3865.7144
$ws.Range("H67").Value = 4179.3
$ws.Range("I67").Value = 4113.7144
$ws.Range("K67").Value = 4113.7144
$ws.Range("M67").Value = -3255.7144
$ws.Range("H92").Value = 565.05
$ws.Range("I92").Value = 557
$ws.Range("J92").Value = 597.25
$ws.Range("K92").Value = 557
$ws.Range("L92").Value = 597.25
$ws.Range("M92").Value = 691
$ws.Range("N92").Value = -3093.25
$ws.Range("H93").Value = 33750
$ws.Range("J93").Value = 33750
$ws.Range("L93").Value = 33750
$ws.Range("N93").Value = -38742
$ws.Range("H100").Value = 1555.5
$ws.Range("I100").Value = 1366.6
$ws.Range("K100").Value = 1366.6
$ws.Range("M100").Value = -825.5999999999999
$ws.Range("H132").Value = 1356.697
$ws.Range("I132").Value = 1357.0385
$ws.Range("K132").Value = 4071.1155
$ws.Range("M132").Value = -1541.1155

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 273103.28
$ws.Range("I32").Value = 333860.38
$ws.Range("J32").Value = 12715.714
$ws.Range("K32").Value = 333860.38
$ws.Range("L32").Value = 12715.714
$ws.Range("M32").Value = -333573.38
$ws.Range("N32").Value = -13289.714
$ws.Range("H45").Value = 42139.64
$ws.Range("I45").Value = 47544.41
$ws.Range("J45").Value = 2504.6667
$ws.Range("K45").Value = 47544.41
$ws.Range("L45").Value = 2504.6667
$ws.Range("M45").Value = -47167.41
$ws.Range("N45").Value = -3258.6667
$ws.Range("H61").Value = 1357330.1
$ws.Range("I61").Value = 40966.32
$ws.Range("K61").Value = 40966.32
$ws.Range("M61").Value = -40754.32
$ws.Range("H74").Value = 324672.25
$ws.Range("I74").Value = 1728.6227
$ws.Range("J74").Value = 1225515
$ws.Range("K74").Value = 1728.6227
$ws.Range("L74").Value = 1225515
$ws.Range("M74").Value = -854.6226999999999
$ws.Range("N74").Value = -1227263
$ws.Range("H77").Value = 324672.25
$ws.Range("I77").Value = 1728.6227
$ws.Range("J77").Value = 1225515
$ws.Range("K77").Value = 8643.113499999999
$ws.Range("L77").Value = 6127575
$ws.Range("M77").Value = -4275.113499999999
$ws.Range("N77").Value = -6136311
$ws.Range("H88").Value = 2378.5334
$ws.Range("I88").Value = 2024.5
$ws.Range("K88").Value = 2024.5
$ws.Range("M88").Value = -1618.5
$ws.Range("H91").Value = 2378.5334
$ws.Range("I91").Value = 2024.5
$ws.Range("K91").Value = 2024.5
$ws.Range("M91").Value = -620.5
$ws.Range("I97").Value = 10005.25
$ws.Range("J97").Value = 3011.889
$ws.Range("K97").Value = 10005.25
$ws.Range("L97").Value = 3011.889
$ws.Range("M97").Value = -9509.25
$ws.Range("N97").Value = -4003.889
$ws.Range("H122").Value = 1579.1111
$ws.Range("I122").Value = 1089
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 3267
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -817
$ws.Range("N122").Value = -21400
$ws.Range("H132").Value = 2844.425
$ws.Range("I132").Value = 1769.3636
$ws.Range("K132").Value = 5308.0908
$ws.Range("M132").Value = -2778.0908
$ws.Range("H135").Value = 64872.223
$ws.Range("J135").Value = 64872.223
$ws.Range("L135").Value = 64872.223
$ws.Range("N135").Value = -75012.223
$ws.Range("H136").Value = 1357330.1
$ws.Range("I136").Value = 40966.32
$ws.Range("K136").Value = 122898.96
$ws.Range("M136").Value = -120348.96

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 942.8125
$ws.Range("I22").Value = 665.6667
$ws.Range("J22").Value = 5100
$ws.Range("K22").Value = 665.6667
$ws.Range("L22").Value = 5100
$ws.Range("M22").Value = -492.6667
$ws.Range("N22").Value = -5446
$ws.Range("H40").Value = 150000
$ws.Range("J40").Value = 150000
$ws.Range("L40").Value = 150000
$ws.Range("N40").Value = -150530
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 2706.2104
$ws.Range("I86").Value = 1173.4286
$ws.Range("K86").Value = 1173.4286
$ws.Range("M86").Value = -50.42859999999996
$ws.Range("H89").Value = 2706.2104
$ws.Range("I89").Value = 1173.4286
$ws.Range("K89").Value = 5867.143
$ws.Range("M89").Value = -251.143
$ws.Range("H94").Value = 1578.0952
$ws.Range("I94").Value = 1660
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 1660
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -1209
$ws.Range("N94").Value = -1702
$ws.Range("H99").Value = 16435.467
$ws.Range("I99").Value = 16579.691
$ws.Range("K99").Value = 16579.691
$ws.Range("M99").Value = -15081.691
$ws.Range("H105").Value = 15979.071
$ws.Range("I105").Value = 13870.8
$ws.Range("J105").Value = 21249.75
$ws.Range("K105").Value = 13870.8
$ws.Range("L105").Value = 21249.75
$ws.Range("M105").Value = -12123.8
$ws.Range("N105").Value = -24743.75
$ws.Range("H107").Value = 19351.691
$ws.Range("I107").Value = 19351.691
$ws.Range("K107").Value = 19351.691
$ws.Range("M107").Value = -17431.691
$ws.Range("H133").Value = 66000
$ws.Range("J133").Value = 66000
$ws.Range("L133").Value = 66000
$ws.Range("N133").Value = -76120
$ws.Range("H138").Value = 109995
$ws.Range("J138").Value = 109995
$ws.Range("L138").Value = 109995
$ws.Range("N138").Value = -120275

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H99").Value = 2001594.9
$ws.Range("I99").Value = 2501514.8
$ws.Range("K99").Value = 2501514.8
$ws.Range("M99").Value = -2500016.8
$ws.Range("H107").Value = 2284
$ws.Range("I107").Value = 2301
$ws.Range("K107").Value = 2301
$ws.Range("M107").Value = -381
$ws.Range("H126").Value = 2001594.9
$ws.Range("I126").Value = 2501514.8
$ws.Range("K126").Value = 7504544.399999999
$ws.Range("M126").Value = -7502074.399999999
$ws.Range("H132").Value = 2291.182
$ws.Range("I132").Value = 2030.3572
$ws.Range("K132").Value = 6091.071599999999
$ws.Range("M132").Value = -3561.071599999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 9507.857
$ws.Range("I62").Value = 9500
$ws.Range("J62").Value = 9509.166999999999
$ws.Range("K62").Value = 28500
$ws.Range("L62").Value = 28527.501
$ws.Range("M62").Value = -27814
$ws.Range("N62").Value = -29899.501
$ws.Range("H65").Value = 9507.857
$ws.Range("I65").Value = 9500
$ws.Range("J65").Value = 9509.166999999999
$ws.Range("K65").Value = 85500
$ws.Range("L65").Value = 85582.503
$ws.Range("M65").Value = -82068
$ws.Range("N65").Value = -92446.503
$ws.Range("H131").Value = 4787285
$ws.Range("I131").Value = 15153549
$ws.Range("J131").Value = 2855.4614
$ws.Range("K131").Value = 45460647
$ws.Range("L131").Value = 8566.3842
$ws.Range("M131").Value = -45455607
$ws.Range("N131").Value = -18646.3842

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1456.2858
$ws.Range("I97").Value = 899.5
$ws.Range("K97").Value = 899.5
$ws.Range("M97").Value = -403.5
$ws.Range("H107").Value = 126024.125
$ws.Range("I107").Value = 500399
$ws.Range("J107").Value = 1232.5
$ws.Range("K107").Value = 500399
$ws.Range("L107").Value = 1232.5
$ws.Range("M107").Value = -498479
$ws.Range("N107").Value = -5072.5
$ws.Range("H122").Value = 3925.7778
$ws.Range("I122").Value = 3476.5715
$ws.Range("K122").Value = 10429.7145
$ws.Range("M122").Value = -7979.7145
$ws.Range("H132").Value = 1834868.9
$ws.Range("I132").Value = 1445.8
$ws.Range("J132").Value = 2853437.2
$ws.Range("K132").Value = 4337.4
$ws.Range("L132").Value = 8560311.600000001
$ws.Range("M132").Value = -1807.4
$ws.Range("N132").Value = -8565371.600000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1295.0769
$ws.Range("I16").Value = 1212.3636
$ws.Range("J16").Value = 1750
$ws.Range("K16").Value = 1212.3636
$ws.Range("L16").Value = 1750
$ws.Range("M16").Value = -1042.3636
$ws.Range("N16").Value = -2090
$ws.Range("H132").Value = 4119.8
$ws.Range("I132").Value = 2867.5557
$ws.Range("K132").Value = 8602.667099999999
$ws.Range("M132").Value = -6072.667099999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1611.3793
$ws.Range("I122").Value = 1496.8572
$ws.Range("J122").Value = 1912
$ws.Range("K122").Value = 4490.571599999999
$ws.Range("L122").Value = 5736
$ws.Range("M122").Value = -2040.571599999999
$ws.Range("N122").Value = -10636
$ws.Range("H132").Value = 2241.9033
$ws.Range("I132").Value = 1946.1177
$ws.Range("J132").Value = 2601.0715
$ws.Range("K132").Value = 5838.3531
$ws.Range("L132").Value = 7803.2145
$ws.Range("M132").Value = -3308.3531
$ws.Range("N132").Value = -12863.2145
